# Generate Report for Handback
# Updates the handback-status workbook with fresh handoff/handback file
# identifiers, hashes and timestamps for the two tracked files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New identifiers
# ---------------------------------------------------------------------
$newId1 = "96764e73-479a-4b60-9e0c-5bab6d3d98fd"
$newId2 = "ffff34319c0c-e71a-4e73-95f4-adba4834507d"

$newHash = "05555d145a175bc9dff08e2daa21167ba2c1cbe8"

$newFileName1 = "$newId1.md"
$newFileName2 = "$newId2.md"
$newPath1 = "e2e\$newId1.md"
$newPath2 = "e2e\$newId2.md"

$newHoDate = "2016-08-15 18:57:11"

$newXlfZhCn = "$newId1.$newHash.zh-cn.xlf"
$newHandoffDateZhCn = "2016-08-15 18:57:01"
$newHandbackDateZhCn = "2016-08-15 18:57:28"

$newXlfDeDe = "$newId1.$newHash.de-de.xlf"
$newHandbackDateDeDe = "2016-08-15 18:57:36"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = $newFileName1
$wsOverview.Range("B2").Value2 = $newPath1
$wsOverview.Range("G2").Value2 = $newHoDate

$wsOverview.Range("A3").Value2 = $newFileName2
$wsOverview.Range("B3").Value2 = $newPath2
$wsOverview.Range("G3").Value2 = $newHoDate

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = $newPath1
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = $newPath2
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = $newFileName1
$wsZhCn.Range("G2").Value2 = $newXlfZhCn
$wsZhCn.Range("H2").Value2 = $newHandoffDateZhCn
$wsZhCn.Range("I2").Value2 = $newFileName1
$wsZhCn.Range("J2").Value2 = $newXlfZhCn
$wsZhCn.Range("K2").Value2 = $newHandbackDateZhCn

$wsZhCn.Range("A3").Value2 = $newFileName2
$wsZhCn.Range("G3").Value2 = $newXlfZhCn
$wsZhCn.Range("H3").Value2 = $newHandoffDateZhCn
$wsZhCn.Range("I3").Value2 = $newFileName2
$wsZhCn.Range("J3").Value2 = $newXlfZhCn
$wsZhCn.Range("K3").Value2 = $newHandbackDateZhCn

foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFileName1
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = $newFileName1
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $newFileName2
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = $newFileName2
    }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = $newFileName1
$wsDeDe.Range("G2").Value2 = $newXlfDeDe
$wsDeDe.Range("H2").Value2 = $newHoDate
$wsDeDe.Range("I2").Value2 = $newFileName1
$wsDeDe.Range("J2").Value2 = $newXlfDeDe
$wsDeDe.Range("K2").Value2 = $newHandbackDateDeDe

$wsDeDe.Range("A3").Value2 = $newFileName2
$wsDeDe.Range("G3").Value2 = $newXlfDeDe
$wsDeDe.Range("H3").Value2 = $newHoDate
$wsDeDe.Range("I3").Value2 = $newFileName2
$wsDeDe.Range("J3").Value2 = $newXlfDeDe
$wsDeDe.Range("K3").Value2 = $newHandbackDateDeDe

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFileName1
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = $newFileName1
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $newFileName2
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = $newFileName2
    }
}
